$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4447565896989806
$ws.Range("C2").Value = 0.6565057718086863
$ws.Range("B3").Value = 1.892803896881188
$ws.Range("C3").Value = 1.48199475587289
$ws.Range("B4").Value = 1.963302924079096
$ws.Range("C4").Value = 2.202644689672959
$ws.Range("B5").Value = 11.03952527224785
$ws.Range("C5").Value = 2.819880456966708
$ws.Range("B6").Value = 11.10963384970866
$ws.Range("C6").Value = 3.489030363414616
$ws.Range("B7").Value = 14.27474525603334
$ws.Range("C7").Value = 4.169497019391097
$ws.Range("B8").Value = 14.37095072489165
$ws.Range("C8").Value = 4.804509904165455
$ws.Range("B9").Value = 18.74205208775246
$ws.Range("C9").Value = 5.411565764825893
$ws.Range("B10").Value = 18.78614874055903
$ws.Range("C10").Value = 6.206170261994902
$ws.Range("B11").Value = 24.94495519855503
$ws.Range("C11").Value = 6.857907327572797
$ws.Range("B12").Value = 24.97006088317273
$ws.Range("C12").Value = 7.560907300024319
$ws.Range("B13").Value = 31.29742046884013
$ws.Range("C13").Value = 8.099090733474986
$ws.Range("B14").Value = 35.93318512432925
$ws.Range("C14").Value = 8.699604280950959
$ws.Range("B15").Value = 41.97090453438198
$ws.Range("C15").Value = 9.415382107315148
$ws.Range("B16").Value = 42.06688104087918
$ws.Range("C16").Value = 10.04440302385245
$ws.Range("B17").Value = 43.54305436741387
$ws.Range("C17").Value = 10.62383611279626
$ws.Range("B18").Value = 53.22130505348912
$ws.Range("C18").Value = 11.33127820246455
$ws.Range("B19").Value = 53.24876112344386
$ws.Range("C19").Value = 11.8834905601153
$ws.Range("B20").Value = 53.31216133828053
$ws.Range("C20").Value = 12.60453721886542
$ws.Range("B21").Value = 53.54367418740043
$ws.Range("C21").Value = 13.27157559785099
$ws.Range("B22").Value = 53.70182378061736
$ws.Range("C22").Value = 13.92005867572347
$ws.Range("B23").Value = 54.34011265091189
$ws.Range("C23").Value = 14.50728062593907
$ws.Range("B24").Value = 54.43735986038732
$ws.Range("C24").Value = 15.11561232594714
$ws.Range("B25").Value = 63.84582133689691
$ws.Range("C25").Value = 15.84623366404085
$ws.Range("B26").Value = 63.86184146201311
$ws.Range("C26").Value = 16.45780705157247
$ws.Range("B27").Value = 64.2621994286214
$ws.Range("C27").Value = 17.22225999458368
$ws.Range("B28").Value = 65.19607471636216
$ws.Range("C28").Value = 17.83300065858448
$ws.Range("B29").Value = 65.24829109583989
$ws.Range("C29").Value = 18.45201313404672
$ws.Range("B30").Value = 67.82466270367233
$ws.Range("C30").Value = 19.33214193278864
$ws.Range("B31").Value = 74.10810647148182
$ws.Range("C31").Value = 20.00339638990841
$ws.Range("B32").Value = 75.35570377662771
$ws.Range("C32").Value = 20.69583483061929
$ws.Range("B33").Value = 76.02580256991875
$ws.Range("C33").Value = 21.30597758286101
$ws.Range("B34").Value = 76.07730442948832
$ws.Range("C34").Value = 21.89811828742358
$ws.Range("B35").Value = 77.96586009814656
$ws.Range("C35").Value = 22.61408567337802
$ws.Range("B36").Value = 80.10168073858485
$ws.Range("C36").Value = 23.20804919705894
$ws.Range("B37").Value = 80.77104262082939
$ws.Range("C37").Value = 23.82593474445101
$ws.Range("B38").Value = 80.83129552967752
$ws.Range("C38").Value = 24.43547901529234
$ws.Range("B39").Value = 81.55346959021148
$ws.Range("C39").Value = 25.08886229050847
$ws.Range("B40").Value = 81.59350169584435
$ws.Range("C40").Value = 25.79238325104663
$ws.Range("B41").Value = 88.91141908178231
$ws.Range("C41").Value = 26.6523679197785
$ws.Range("B42").Value = 90.17469245606924
$ws.Range("C42").Value = 27.20179538036021
$ws.Range("B43").Value = 90.24982650559828
$ws.Range("C43").Value = 27.8551061257309
$ws.Range("B44").Value = 92.26533388883108
$ws.Range("C44").Value = 28.43591586257021
$ws.Range("B45").Value = 92.34116824077697
$ws.Range("C45").Value = 29.18459100326947
$ws.Range("B46").Value = 92.50596827054754
$ws.Range("C46").Value = 29.91138212360896
$ws.Range("B47").Value = 99.80538993828868
$ws.Range("C47").Value = 30.59144664013221
$ws.Range("B48").Value = 99.87672250723311
$ws.Range("C48").Value = 31.21469843703736
$ws.Range("B49").Value = 99.91903307001185
$ws.Range("C49").Value = 31.8637692663708
